# chelyabinsk.xlsx update
# Commit message: "Add 2 sites spb, fix omsk, del 3 sites"
#
# Net data changes on the single worksheet ("Sheet"):
#  1. Lada "Largus New" (row 161): the saturn2.ru price/link (J:K) is removed.
#  2. UAZ "3909 Бортовой грузовик" (id 695, row 229) is removed entirely -
#     all rows below it shift up by one.
#  3. XCite "X-Cross 8" (id 787), which ends up on row 244 after the shift,
#     loses its saturn2.ru price/link (J:K).
#
# The used-range dimension (A1:Q245 -> A1:Q244) updates automatically as a
# consequence of deleting the whole row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the saturn2.ru site columns for "Largus New" (row 161).
$ws.Range("J161:K161").ClearContents()

# 2) Delete the UAZ "3909 Бортовой грузовик" row outright (shifts rows up).
$ws.Range("A229").EntireRow.Delete()

# 3) After the shift, "X-Cross 8" (formerly row 245) is now row 244;
#    drop its saturn2.ru site columns.
$ws.Range("J244:K244").ClearContents()
